$d = $word.ActiveDocument

# --- Change 1: the lone "b)" list item right before the "d) precision..." item
# becomes "a)". There are multiple stand-alone "b)" paragraphs in the document,
# so locate the specific one that is immediately followed by the paragraph
# starting with "d)" and containing "precision = 0.85".
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $txt = $para.Range.Text
    if ($txt.TrimEnd("`r`a") -eq "b)") {
        if ($i -lt $d.Paragraphs.Count) {
            $nextTxt = $d.Paragraphs($i + 1).Range.Text
            if ($nextTxt -like "d)*precision*") {
                $target = $para
                break
            }
        }
    }
}

if ($target -ne $null) {
    $rng = $target.Range
    $rng.Find.Execute("b)", $true, $false, $false, $false, $false, $true, 1, $false, "a)", 2)
}

# --- Change 2: merge the standalone " " run with the following
# "precision = 0.85 = VP / (VP + FP) " run into a single run
# (" precision = 0.85 = VP / (VP + FP) "), without touching the preceding
# "d)" run. We restrict the Find range to start exactly after "d) " so the
# replace only coalesces the "precision..." run with the immediately
# preceding space run, leaving "d)" as its own run.
$searchText = "precision = 0.85 = VP / (VP + FP) "
$full = $d.Content
$full.Find.ClearFormatting()
$found = $full.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $full.Start
    $end = $full.End
    $scoped = $d.Range($start, $end)
    $scoped.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $searchText, 2)
}
